# Commit: "Kontakt bearbeiten und hinzufügen"
#
# Appends two new logged time entries to the "Zeitliste Maturaprojekt"
# tracking sheet (Tabelle1):
#   Row 37: 2020-01-15 (serial 43845) - "Kontakt hinzufügen, bearbeiten"        - 1   hour
#   Row 38: 2020-01-16 (serial 43846) - "Organisatorisches; Kontakt hinzufügen, bearbeiten" - 3.5 hours
#
# The F3 (=SUM(C4:C100)) and F4 (=180-F3) totals recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the date cell format (style "2", a date number format) from the
# last existing data row onto the two new rows, so A37/A38 render as dates
# just like A4:A36 without introducing a brand-new style/numFmt entry.
$ws.Range("A36").Copy()
$ws.Range("A37:A38").PasteSpecial(-4122)

# Row 37
$ws.Range("A37").Value = 43845
$ws.Range("B37").Value = "Kontakt hinzufügen, bearbeiten"
$ws.Range("C37").Value = 1

# Row 38
$ws.Range("A38").Value = 43846
$ws.Range("B38").Value = "Organisatorisches; Kontakt hinzufügen, bearbeiten"
$ws.Range("C38").Value = 3.5

# Move the active selection to follow the newly added data, matching the
# author's cursor position after entering the last row.
[void]$ws.Range("B39").Select()
